# Natmi following Dr Hou advice
# Rewrites the LR-pair table (rows 2-21, cols A:T) on the active sheet
# with the updated sending/target cluster pairing and recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = "ECs"
$data[0,1] = "Tgfb1"
$data[0,2] = "Sdc2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 50.950239
$data[0,7] = 152.850717
$data[0,8] = 0.1520006117784607
$data[0,9] = 0.1540898474582185
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.511114
$data[0,13] = 4.533342
$data[0,14] = 0.02241848840500565
$data[0,15] = 0.02483957450832654
$data[0,16] = 76.99161945624601
$data[0,17] = 692.9245751062141
$data[0,18] = 0.003407623952709187
$data[0,19] = 0.003827526246915089
$data[1,0] = "ECs"
$data[1,1] = "Tgfb1"
$data[1,2] = "Sdc2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 50.950239
$data[1,7] = 152.850717
$data[1,8] = 0.1520006117784607
$data[1,9] = 0.1540898474582185
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 46.15376066666666
$data[1,13] = 138.461282
$data[1,14] = 0.6847250097299556
$data[1,15] = 0.7586719313825014
$data[1,16] = 2351.545136715466
$data[1,17] = 21163.90623043919
$data[1,18] = 0.1040786203789657
$data[1,19] = 0.1169036421775616
$data[2,0] = "ECs"
$data[2,1] = "Tgfb1"
$data[2,2] = "Sdc2"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 50.950239
$data[2,7] = 152.850717
$data[2,8] = 0.1520006117784607
$data[2,9] = 0.1540898474582185
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.030318
$data[2,13] = 0.09095399999999999
$data[2,14] = 0.0004497898447522563
$data[2,15] = 0.000498364928088446
$data[2,16] = 1.544709346002
$data[2,17] = 13.902384114018
$data[2,18] = [double]"6.836833157408184E-05"
$data[2,19] = [double]"7.679297574767467E-05"
$data[3,0] = "ECs"
$data[3,1] = "Tgfb1"
$data[3,2] = "Sdc2"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 50.950239
$data[3,7] = 152.850717
$data[3,8] = 0.1520006117784607
$data[3,9] = 0.1540898474582185
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 19.7096195
$data[3,13] = 39.419239
$data[3,14] = 0.2924067120202865
$data[3,15] = 0.2159901291810835
$data[3,16] = 1004.209824124061
$data[3,17] = 6025.258944744364
$data[3,18] = 0.04444599911521174
$data[3,19] = 0.03328188605799406
$data[4,0] = "FAPs"
$data[4,1] = "Tgfb1"
$data[4,2] = "Sdc2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 15.19715733333333
$data[4,7] = 45.591472
$data[4,8] = 0.04533790728558088
$data[4,9] = 0.0459610730244441
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.511114
$data[4,13] = 4.533342
$data[4,14] = 0.02241848840500565
$data[4,15] = 0.02483957450832654
$data[4,16] = 22.96463720660266
$data[4,17] = 206.681734859424
$data[4,18] = 0.001016407348789016
$data[4,19] = 0.001141653497873316
$data[5,0] = "FAPs"
$data[5,1] = "Tgfb1"
$data[5,2] = "Sdc2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 15.19715733333333
$data[5,7] = 45.591472
$data[5,8] = 0.04533790728558088
$data[5,9] = 0.0459610730244441
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 46.15376066666666
$data[5,13] = 138.461282
$data[5,14] = 0.6847250097299556
$data[5,15] = 0.7586719313825014
$data[5,16] = 701.4059623763447
$data[5,17] = 6312.653661387103
$data[5,18] = 0.03104399900725519
$data[5,19] = 0.03486937603986719
$data[6,0] = "FAPs"
$data[6,1] = "Tgfb1"
$data[6,2] = "Sdc2"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 15.19715733333333
$data[6,7] = 45.591472
$data[6,8] = 0.04533790728558088
$data[6,9] = 0.0459610730244441
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.030318
$data[6,13] = 0.09095399999999999
$data[6,14] = 0.0004497898447522563
$data[6,15] = 0.000498364928088446
$data[6,16] = 0.4607474160319999
$data[6,17] = 4.146726744287999
$data[6,18] = [double]"2.039253027937361E-05"
$data[6,19] = [double]"2.29053868526949E-05"
$data[7,0] = "FAPs"
$data[7,1] = "Tgfb1"
$data[7,2] = "Sdc2"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 15.19715733333333
$data[7,7] = 45.591472
$data[7,8] = 0.04533790728558088
$data[7,9] = 0.0459610730244441
$data[7,10] = 2
$data[7,11] = 1
$data[7,12] = 19.7096195
$data[7,13] = 39.419239
$data[7,14] = 0.2924067120202865
$data[7,15] = 0.2159901291810835
$data[7,16] = 299.5301885216347
$data[7,17] = 1797.181131129808
$data[7,18] = 0.0132571083992573
$data[7,19] = 0.009927138099850893
$data[8,0] = "M1"
$data[8,1] = "Tgfb1"
$data[8,2] = "Sdc2"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 124.2523
$data[8,7] = 372.7569
$data[8,8] = 0.3706837491945981
$data[8,9] = 0.375778766284743
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 1.511114
$data[8,13] = 4.533342
$data[8,14] = 0.02241848840500565
$data[8,15] = 0.02483957450832654
$data[8,16] = 187.7593900622
$data[8,17] = 1689.8345105598
$data[8,18] = 0.00831016933324312
$data[8,19] = 0.0093341846637769
$data[9,0] = "M1"
$data[9,1] = "Tgfb1"
$data[9,2] = "Sdc2"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 124.2523
$data[9,7] = 372.7569
$data[9,8] = 0.3706837491945981
$data[9,9] = 0.375778766284743
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 46.15376066666666
$data[9,13] = 138.461282
$data[9,14] = 0.6847250097299556
$data[9,15] = 0.7586719313825014
$data[9,16] = 5734.710916482866
$data[9,17] = 51612.39824834579
$data[9,18] = 0.2538164337740076
$data[9,19] = 0.2850928023897796
$data[10,0] = "M1"
$data[10,1] = "Tgfb1"
$data[10,2] = "Sdc2"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 124.2523
$data[10,7] = 372.7569
$data[10,8] = 0.3706837491945981
$data[10,9] = 0.375778766284743
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.030318
$data[10,13] = 0.09095399999999999
$data[10,14] = 0.0004497898447522563
$data[10,15] = 0.000498364928088446
$data[10,16] = 3.767081231399999
$data[10,17] = 33.9037310826
$data[10,18] = 0.0001667297860024226
$data[10,19] = 0.0001872749578366609
$data[11,0] = "M1"
$data[11,1] = "Tgfb1"
$data[11,2] = "Sdc2"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 124.2523
$data[11,7] = 372.7569
$data[11,8] = 0.3706837491945981
$data[11,9] = 0.375778766284743
$data[11,10] = 2
$data[11,11] = 1
$data[11,12] = 19.7096195
$data[11,13] = 39.419239
$data[11,14] = 0.2924067120202865
$data[11,15] = 0.2159901291810835
$data[11,16] = 2448.96555499985
$data[11,17] = 14693.7933299991
$data[11,18] = 0.108390416301345
$data[11,19] = 0.08116450427334984
$data[12,0] = "M2"
$data[12,1] = "Tgfb1"
$data[12,2] = "Sdc2"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 131.1635103333333
$data[12,7] = 393.490531
$data[12,8] = 0.3913020665845575
$data[12,9] = 0.3966804807205673
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 1.511114
$data[12,13] = 4.533342
$data[12,14] = 0.02241848840500565
$data[12,15] = 0.02483957450832654
$data[12,16] = 198.2030167538447
$data[12,17] = 1783.827150784602
$data[12,18] = 0.008772400842580651
$data[12,19] = 0.009853374356857321
$data[13,0] = "M2"
$data[13,1] = "Tgfb1"
$data[13,2] = "Sdc2"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 131.1635103333333
$data[13,7] = 393.490531
$data[13,8] = 0.3913020665845575
$data[13,9] = 0.3966804807205673
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 46.15376066666666
$data[13,13] = 138.461282
$data[13,14] = 0.6847250097299556
$data[13,15] = 0.7586719313825014
$data[13,16] = 6053.689264124527
$data[13,17] = 54483.20337712074
$data[13,18] = 0.2679343113494629
$data[13,19] = 0.300950346450012
$data[14,0] = "M2"
$data[14,1] = "Tgfb1"
$data[14,2] = "Sdc2"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 131.1635103333333
$data[14,7] = 393.490531
$data[14,8] = 0.3913020665845575
$data[14,9] = 0.3966804807205673
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.030318
$data[14,13] = 0.09095399999999999
$data[14,14] = 0.0004497898447522563
$data[14,15] = 0.000498364928088446
$data[14,16] = 3.976615306286
$data[14,17] = 35.789537756574
$data[14,18] = 0.0001760036957803052
$data[14,19] = 0.0001976916392483957
$data[15,0] = "M2"
$data[15,1] = "Tgfb1"
$data[15,2] = "Sdc2"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 131.1635103333333
$data[15,7] = 393.490531
$data[15,8] = 0.3913020665845575
$data[15,9] = 0.3966804807205673
$data[15,10] = 2
$data[15,11] = 1
$data[15,12] = 19.7096195
$data[15,13] = 39.419239
$data[15,14] = 0.2924067120202865
$data[15,15] = 0.2159901291810835
$data[15,16] = 2585.182880954319
$data[15,17] = 15511.09728572591
$data[15,18] = 0.1144193506967337
$data[15,19] = 0.08567906827444964
$data[16,0] = "sCs"
$data[16,1] = "Tgfb1"
$data[16,2] = "Sdc2"
$data[16,3] = "ECs"
$data[16,4] = 2
$data[16,5] = 1
$data[16,6] = 13.634385
$data[16,7] = 27.26877
$data[16,8] = 0.04067566515680266
$data[16,9] = 0.02748983251202704
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 1.511114
$data[16,13] = 4.533342
$data[16,14] = 0.02241848840500565
$data[16,15] = 0.02483957450832654
$data[16,16] = 20.60311005489
$data[16,17] = 123.61866032934
$data[16,18] = 0.0009118869276836727
$data[16,19] = 0.000682835742903913
$data[17,0] = "sCs"
$data[17,1] = "Tgfb1"
$data[17,2] = "Sdc2"
$data[17,3] = "FAPs"
$data[17,4] = 2
$data[17,5] = 1
$data[17,6] = 13.634385
$data[17,7] = 27.26877
$data[17,8] = 0.04067566515680266
$data[17,9] = 0.02748983251202704
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 46.15376066666666
$data[17,13] = 138.461282
$data[17,14] = 0.6847250097299556
$data[17,15] = 0.7586719313825014
$data[17,16] = 629.27814212719
$data[17,17] = 3775.66885276314
$data[17,18] = 0.02785164522026412
$data[17,19] = 0.02085576432528104
$data[18,0] = "sCs"
$data[18,1] = "Tgfb1"
$data[18,2] = "Sdc2"
$data[18,3] = "M2"
$data[18,4] = 2
$data[18,5] = 1
$data[18,6] = 13.634385
$data[18,7] = 27.26877
$data[18,8] = 0.04067566515680266
$data[18,9] = 0.02748983251202704
$data[18,10] = 1
$data[18,11] = 0.3333333333333333
$data[18,12] = 0.030318
$data[18,13] = 0.09095399999999999
$data[18,14] = 0.0004497898447522563
$data[18,15] = 0.000498364928088446
$data[18,16] = 0.41336728443
$data[18,17] = 2.48020370658
$data[18,18] = [double]"1.829550111607303E-05"
$data[18,19] = [double]"1.369996840301978E-05"
$data[19,0] = "sCs"
$data[19,1] = "Tgfb1"
$data[19,2] = "Sdc2"
$data[19,3] = "sCs"
$data[19,4] = 2
$data[19,5] = 1
$data[19,6] = 13.634385
$data[19,7] = 27.26877
$data[19,8] = 0.04067566515680266
$data[19,9] = 0.02748983251202704
$data[19,10] = 2
$data[19,11] = 1
$data[19,12] = 19.7096195
$data[19,13] = 39.419239
$data[19,14] = 0.2924067120202865
$data[19,15] = 0.2159901291810835
$data[19,16] = 268.7285404665075
$data[19,17] = 1074.91416186603
$data[19,18] = 0.0118938375077388
$data[19,19] = 0.00593753247543907

$ws.Range("A2:T21").Value = $data
